# Scoreboard.xlsx update
# 1) Rename the "Sofies Disipler" team (row 18 of ScoreF, partner of "Hulda og Kristiane")
#    to "In it for the party!"
# 2) Fill in the Minute/Second/Rep-style results (columns D, E, F) for rows 2-26 on
#    the ScoreF sheet, which previously only had the Name columns (A, B) filled in.
# 3) Make ScoreF the active sheet/tab, with F22 selected (ScoreM keeps its F23 selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScoreF")

# Rename the team in column B, row 18
$ws.Range("B18").Value = "In it for the party!"

# Fill in D (minutes), E (seconds), F (reps) results for each team row
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 39
$ws.Cells.Item(2, 6).Value = 380

$ws.Cells.Item(3, 4).Value = 6
$ws.Cells.Item(3, 5).Value = 50
$ws.Cells.Item(3, 6).Value = 380

$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 17
$ws.Cells.Item(4, 6).Value = 380

$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 22
$ws.Cells.Item(5, 6).Value = 380

$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 37
$ws.Cells.Item(6, 6).Value = 380

$ws.Cells.Item(7, 4).Value = 6
$ws.Cells.Item(7, 5).Value = 19
$ws.Cells.Item(7, 6).Value = 380

$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 39
$ws.Cells.Item(8, 6).Value = 380

$ws.Cells.Item(9, 4).Value = 6
$ws.Cells.Item(9, 5).Value = 23
$ws.Cells.Item(9, 6).Value = 380

$ws.Cells.Item(10, 4).Value = 6
$ws.Cells.Item(10, 5).Value = 12
$ws.Cells.Item(10, 6).Value = 380

$ws.Cells.Item(11, 4).Value = 7
$ws.Cells.Item(11, 5).Value = 35
$ws.Cells.Item(11, 6).Value = 380

$ws.Cells.Item(12, 4).Value = 6
$ws.Cells.Item(12, 5).Value = 55
$ws.Cells.Item(12, 6).Value = 380

$ws.Cells.Item(13, 4).Value = 7
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 380

$ws.Cells.Item(14, 4).Value = 8
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 355

$ws.Cells.Item(15, 4).Value = 7
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 380

$ws.Cells.Item(16, 4).Value = 5
$ws.Cells.Item(16, 5).Value = 31
$ws.Cells.Item(16, 6).Value = 380

$ws.Cells.Item(17, 4).Value = 7
$ws.Cells.Item(17, 5).Value = 52
$ws.Cells.Item(17, 6).Value = 380

$ws.Cells.Item(18, 4).Value = 8
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 330

$ws.Cells.Item(19, 4).Value = 5
$ws.Cells.Item(19, 5).Value = 18
$ws.Cells.Item(19, 6).Value = 330

$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = 20
$ws.Cells.Item(20, 6).Value = 380

$ws.Cells.Item(21, 4).Value = 6
$ws.Cells.Item(21, 5).Value = 27
$ws.Cells.Item(21, 6).Value = 380

$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(22, 5).Value = 35
$ws.Cells.Item(22, 6).Value = 380

$ws.Cells.Item(23, 4).Value = 6
$ws.Cells.Item(23, 5).Value = 14
$ws.Cells.Item(23, 6).Value = 380

$ws.Cells.Item(24, 4).Value = 6
$ws.Cells.Item(24, 5).Value = 42
$ws.Cells.Item(24, 6).Value = 380

$ws.Cells.Item(25, 4).Value = 7
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = 380

$ws.Cells.Item(26, 4).Value = 6
$ws.Cells.Item(26, 5).Value = 36
$ws.Cells.Item(26, 6).Value = 380

# Activate ScoreF and update its selection (the commit made ScoreF the visible/active tab)
$ws.Activate() | Out-Null
$ws.Range("F22").Select() | Out-Null
